$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.541.94"
$ws.Range("E2").Value = "  -3.85%  "
$ws.Range("D3").Value = "2.371.39"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.54"
$ws.Range("E5").Value = "  -5.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.37"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  -3.19%  "
$ws.Range("D9").Value = "2.390.67"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0961"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.320"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("E13").Value = "  -10.05%  "
$ws.Range("D14").Value = "2.796.10"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("D15").Value = "56.350.79"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.58"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "2.384.91"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.20"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.96"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.04"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.74"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "2.495.32"
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.150"
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.370"
$ws.Range("E28").Value = "  -9.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.31"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.66"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.13"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.76"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  -4.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.79"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").Value = "  -5.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.779"
$ws.Range("E42").Value = "  -6.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.11"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.35"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "254.59"
$ws.Range("E46").Value = "  -6.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.568"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0899"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0487"
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.81"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("E51").Value = "  -4.53%  "

Write-Host "Applied all cryptos list updates"
